$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E23").Value = "https://ajc.maths.uq.edu.au/pdf/20/ocr-ajc-v20-p19.pdf"
$ws.Range("B23").Value = "Landau's Theorem revisited"
$ws.Range("C23").Value = "Jerrold Griggs, K. B. Reid"
$ws.Range("F23").Value = "Round Robins"
$ws.Range("G23").Value = "?"
$ws.Range("H23").Value = "nice proof"
$ws.Range("D23").Value = 1999

$ws.Range("H24").Select()
